$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers for new columns I and J ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header styling used by the existing header row (copy format from H1,
# which carries the bold/bordered/centered "header" style).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- Data for columns I (I0) and J (IF), rows 2-56 ---
$i0 = @(8,8,7,6,7,9,9,8,8,7,9,7,6,6,9,8,7,7,7,6,7,6,6,1,6,2,8,1,7,6,7,7,7,7,8,8,6,10,7,8,7,6,6,6,6,7,6,8,6,5,9,5,8,6,3)
$if = @(8,8,7,6,7,9,9,9,8,8,10,7,6,6,9,9,7,7,7,6,7,6,6,1,7,3,8,1,7,6,8,7,7,7,8,8,6,10,7,8,7,6,6,6,6,7,6,8,6,5,9,6,8,7,3)

for ($idx = 0; $idx -lt $i0.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $i0[$idx]
    $ws.Cells.Item($row, 10).Value = $if[$idx]
}
